$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new "2023" column (S) to the table, reusing the same cell
# formatting as the existing "2022" column (R) for each row. Inserting
# the copied range (rather than a plain paste) is what makes the new
# cells pick up the same style indices Excel itself would reuse.
$src = $ws.Range("R3:R14")
$src.Copy()
$dst = $ws.Range("S3:S14")
$dst.Insert(-4161)   # xlShiftToRight

# Fill in the 2023 values.
$ws.Range("S3").Value = 2023
$ws.Range("S4").Value = 89.1
$ws.Range("S5").Value = 89.1
$ws.Range("S6").Value = 1895
$ws.Range("S7").Value = 1759
$ws.Range("S8").Value = 683.8
$ws.Range("S9").Value = 40.7
$ws.Range("S10").Value = 14.7
$ws.Range("S11").Value = 48.4
$ws.Range("S12").Value = 0.2
$ws.Range("S13").Value = 40.5
$ws.Range("S14").Value = "_"

# Leave the selection where it ended up after entering the data.
$ws.Range("O22").Select()
